# CMS_Avionics_Channels.xlsx — sync channel config with config.hpp
#
# The "channels" sheet lists avionics channel IDs in column B. Two
# adjacent-row ID pairs were swapped relative to config.hpp, so fix them:
#   - FU_LOWER_SETP (row 12)    / OX_UPPER_SETP (row 13)
#   - FU_LOWER_REDLINE (row 16) / OX_UPPER_REDLINE (row 17)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("channels")

$ws.Range("B12").Value = 16
$ws.Range("B13").Value = 15
$ws.Range("B16").Value = 20
$ws.Range("B17").Value = 19

# Leave the selection where the edit finished, as Excel would.
$ws.Range("B13").Select()
